$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename first sheet
$ws1.Name = "Sheet1"

# Remove special style from A2/A3 (set to default style "Normal")
$ws1.Range("A2").Style = "Normal"
$ws1.Range("A3").Style = "Normal"

# Update existing row 2 (aapl)
$ws1.Cells.Item(2, 4).Value = 272.1900024414062

# Update existing row 3 (goog)
$ws1.Cells.Item(3, 2).Value = 14296
$ws1.Cells.Item(3, 3).Value = 204.2062094332519
$ws1.Cells.Item(3, 4).Value = 303.75
$ws1.Cells.Item(3, 5).Value = 50.69571676610289

# Add new rows 4-6
$ws1.Cells.Item(4, 1).Value = "amzn"
$ws1.Cells.Item(4, 2).Value = 129
$ws1.Cells.Item(4, 3).Value = 237.4199981689453
$ws1.Cells.Item(4, 4).Value = 226.7599945068359
$ws1.Cells.Item(4, 5).Value = -6.802288779815519

$ws1.Cells.Item(5, 1).Value = "celh"
$ws1.Cells.Item(5, 2).Value = 120
$ws1.Cells.Item(5, 3).Value = 23.76000022888184
$ws1.Cells.Item(5, 4).Value = 41.66999816894531
$ws1.Cells.Item(5, 5).Value = 71.8434403945129

$ws1.Cells.Item(6, 1).Value = "rivn"
$ws1.Cells.Item(6, 2).Value = 150
$ws1.Cells.Item(6, 3).Value = 13.8100004196167
$ws1.Cells.Item(6, 4).Value = 20.28000068664551
$ws1.Cells.Item(6, 5).Value = 46.8501091269944

# Update Summary sheet row 2
$ws2.Cells.Item(2, 2).Value = 5
$ws2.Cells.Item(2, 3).Value = 3234979.363583846
$ws2.Cells.Item(2, 4).Value = 15929
$ws2.Cells.Item(2, 5).Value = 203.0874106085659
$ws2.Cells.Item(2, 6).Value = 22.44052187499769
